$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the data block (rows 540-543), pushing the
# existing rows 540-577 down to 544-581.
$ws.Rows("540:543").Insert()

# Row 540
$ws.Range("A540").Value = 1
$ws.Range("B540").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C540").Value = "Arica y Parinacota"
$ws.Range("D540").Value = 44610
$ws.Range("E540").Value = 15
$ws.Range("F540").Value = 100112024
$ws.Range("G540").Value = "Choclo"
$ws.Range("H540").Value = "Lluteño"
$ws.Range("I540").Value = "Primera"
$ws.Range("J540").Value = 50
$ws.Range("K540").Value = 23000
$ws.Range("L540").Value = 24000
$ws.Range("M540").Value = 23500
$ws.Range("N540").Value = "$/saco 50 unidades"
$ws.Range("O540").Value = "Región de Arica y Parinacota"
$ws.Range("P540").Value = 470
$ws.Range("Q540").Value = 50
$ws.Range("R540").Value = "Hortaliza"

# Row 541
$ws.Range("A541").Value = 1
$ws.Range("B541").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C541").Value = "Arica y Parinacota"
$ws.Range("D541").Value = 44610
$ws.Range("E541").Value = 15
$ws.Range("F541").Value = 100112024
$ws.Range("G541").Value = "Choclo"
$ws.Range("H541").Value = "Lluteño"
$ws.Range("I541").Value = "Segunda"
$ws.Range("J541").Value = 60
$ws.Range("K541").Value = 21000
$ws.Range("L541").Value = 22000
$ws.Range("M541").Value = 21500
$ws.Range("N541").Value = "$/saco 75 unidades"
$ws.Range("O541").Value = "Región de Arica y Parinacota"
$ws.Range("P541").Value = 287
$ws.Range("Q541").Value = 75
$ws.Range("R541").Value = "Hortaliza"

# Row 542
$ws.Range("A542").Value = 1
$ws.Range("B542").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C542").Value = "Arica y Parinacota"
$ws.Range("D542").Value = 44610
$ws.Range("E542").Value = 15
$ws.Range("F542").Value = 100112024
$ws.Range("G542").Value = "Choclo"
$ws.Range("H542").Value = "Lluteño"
$ws.Range("I542").Value = "Tercera"
$ws.Range("J542").Value = 60
$ws.Range("K542").Value = 19000
$ws.Range("L542").Value = 20000
$ws.Range("M542").Value = 19500
$ws.Range("N542").Value = "$/saco 100 unidades"
$ws.Range("O542").Value = "Región de Arica y Parinacota"
$ws.Range("P542").Value = 195
$ws.Range("Q542").Value = 100
$ws.Range("R542").Value = "Hortaliza"

# Row 543
$ws.Range("A543").Value = 1
$ws.Range("B543").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C543").Value = "Arica y Parinacota"
$ws.Range("D543").Value = 44610
$ws.Range("E543").Value = 15
$ws.Range("F543").Value = 100112024
$ws.Range("G543").Value = "Choclo"
$ws.Range("H543").Value = "Sin especificar"
$ws.Range("I543").Value = "Primera"
$ws.Range("J543").Value = 40
$ws.Range("K543").Value = 28000
$ws.Range("L543").Value = 30000
$ws.Range("M543").Value = 29000
$ws.Range("N543").Value = "$/saco 100 unidades"
$ws.Range("O543").Value = "Valle de Camiña"
$ws.Range("P543").Value = 290
$ws.Range("Q543").Value = 100
$ws.Range("R543").Value = "Hortaliza"
